$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Z1: new date header cell, stored as literal text (matches source style of
#     the existing date-header row, which was authored by openpyxl as inlineStr,
#     not an Excel date serial). Pre-formatting the cell as Text keeps Excel from
#     auto-converting the "dd-mm-yyyy"-looking string into a date value.
$z1 = $ws.Range("Z1")
$z1.NumberFormat = "@"
$z1.Value = "11-10-2020"

# Match the visual formatting of the neighboring header cells (bold, centered,
# thin box border) like Y1 uses.
$z1.Font.Bold = $true
$z1.HorizontalAlignment = -4108   # xlCenter
$z1.VerticalAlignment = -4160     # xlTop
$z1.Borders.LineStyle = 1         # xlContinuous
$z1.Borders.Weight = 2            # xlThin

# --- Z2:Z36: new daily deceased-count column (values sourced from the updated
#     time series for 11-10-2020).
$ws.Cells.Item(2, 26).Value = 55
$ws.Cells.Item(3, 26).Value = 6194
$ws.Cells.Item(4, 26).Value = 23
$ws.Cells.Item(5, 26).Value = 811
$ws.Cells.Item(6, 26).Value = 944
$ws.Cells.Item(7, 26).Value = 190
$ws.Cells.Item(8, 26).Value = 1235
$ws.Cells.Item(9, 26).Value = 2
$ws.Cells.Item(10, 26).Value = 5740
$ws.Cells.Item(11, 26).Value = 499
$ws.Cells.Item(12, 26).Value = 3557
$ws.Cells.Item(13, 26).Value = 1572
$ws.Cells.Item(14, 26).Value = 248
$ws.Cells.Item(15, 26).Value = 1313
$ws.Cells.Item(16, 26).Value = 784
$ws.Cells.Item(17, 26).Value = 9891
$ws.Cells.Item(18, 26).Value = 978
$ws.Cells.Item(19, 26).Value = 64
$ws.Cells.Item(20, 26).Value = 2599
$ws.Cells.Item(21, 26).Value = 40040
$ws.Cells.Item(22, 26).Value = 88
$ws.Cells.Item(23, 26).Value = 62
$ws.Cells.Item(24, 26).Value = 0
$ws.Cells.Item(25, 26).Value = 17
$ws.Cells.Item(26, 26).Value = 1006
$ws.Cells.Item(27, 26).Value = 559
$ws.Cells.Item(28, 26).Value = 3798
$ws.Cells.Item(29, 26).Value = 1636
$ws.Cells.Item(30, 26).Value = 55
$ws.Cells.Item(31, 26).Value = 10187
$ws.Cells.Item(32, 26).Value = 1222
$ws.Cells.Item(33, 26).Value = 315
$ws.Cells.Item(34, 26).Value = 734
$ws.Cells.Item(35, 26).Value = 6353
$ws.Cells.Item(36, 26).Value = 5563
